$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-33)
# from serial date 45175 (2023-09-06) to 45177 (2023-09-08)
for ($row = 2; $row -le 33; $row++) {
    $ws.Cells.Item($row, 3).Value = 45177
}
